# Add the new "AddEmergencyContactTest" worksheet right after the existing
# "InvalidCredentialTest" sheet (becomes sheetId=2 / rId2, and the newly
# inserted sheet automatically becomes the active/selected tab).
$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws  = $wb.Worksheets.Add($null, $ws1)
$ws.Name = "AddEmergencyContactTest"

# F2/G2 must end up as real numbers (34343 / 343434) even though the rest of
# the sheet uses a Text ("@") number format. Writing them while the cell
# format is still the default General keeps them numeric; applying the Text
# format afterwards (next step) does not retroactively convert them.
$ws.Cells.Item(2, 6).Value = 34343
$ws.Cells.Item(2, 7).Value = 343434

# Apply the Text number format ("@", numFmtId 49) to the whole used range.
# Cells that still get a plain numeric-looking string written into them
# after this point (e.g. "6767") are stored as text/shared-strings, matching
# the source data which mixes genuine numbers with zero-padded/text phone
# numbers.
$rng = $ws.Range("A1:H4")
$rng.NumberFormat = "@"

# Row 1 - header labels.
# NOTE: columns are populated in the order A,B,C,D,G,F,E,H (not strictly
# left-to-right) so that the workbook's shared-string table is built up in
# the same order as in the authored workbook.
$ws.Cells.Item(1, 1).Value = "Username"
$ws.Cells.Item(1, 2).Value = "Password"
$ws.Cells.Item(1, 3).Value = "Contact Name"
$ws.Cells.Item(1, 4).Value = "Relationship"
$ws.Cells.Item(1, 7).Value = "Home Telephone"
$ws.Cells.Item(1, 6).Value = "Mobile Telephone"
$ws.Cells.Item(1, 5).Value = "Home Phone"
$ws.Cells.Item(1, 8).Value = "Expected Value"

# Row 2 (F2/G2 were already written above as numbers).
$ws.Cells.Item(2, 1).Value = "Admin"
$ws.Cells.Item(2, 2).Value = "admin123"
$ws.Cells.Item(2, 3).Value = "Deep"
$ws.Cells.Item(2, 4).Value = "Sister"
$ws.Cells.Item(2, 5).Value = "545454"
$ws.Cells.Item(2, 8).Value = "Deep;545454"

# Row 3
$ws.Cells.Item(3, 1).Value = "Admin"
$ws.Cells.Item(3, 2).Value = "admin123"
$ws.Cells.Item(3, 3).Value = "Faha"
$ws.Cells.Item(3, 4).Value = "Father"
$ws.Cells.Item(3, 5).Value = "676767"
$ws.Cells.Item(3, 6).Value = "6767"
$ws.Cells.Item(3, 7).Value = "6767"
$ws.Cells.Item(3, 8).Value = "Faha;6767;6767"

# Row 4
$ws.Cells.Item(4, 1).Value = "Admin"
$ws.Cells.Item(4, 2).Value = "admin123"
$ws.Cells.Item(4, 3).Value = "Faha12"
$ws.Cells.Item(4, 4).Value = "Father"
$ws.Cells.Item(4, 5).Value = "676767"
$ws.Cells.Item(4, 6).Value = "6767"
$ws.Cells.Item(4, 7).Value = "6767"
$ws.Cells.Item(4, 8).Value = "Faha12;6767;6767"

# Match the authored selection on the new sheet.
$null = $ws.Range("G7").Select()
